# Fix the typo "scratching" -> "scraping" in the "second step" bullet.
# "scratching" = "scra" + "tch" + "ing"; the middle "tch" becomes "p",
# giving "scra" + "p" + "ing" = "scraping". The rest of the paragraph
# (e.g. "...faisant du ") is left untouched.
$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("scratching", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $wordStart = $rng.Start

    $tchStart = $wordStart + 4
    $tchEnd = $wordStart + 7

    $rTch = $d.Range($tchStart, $tchEnd)
    $rTch.Text = "p"
} else {
    # Fallback: plain whole-word replace in case the locate-and-patch
    # above could not find the expected text.
    $d.Content.Find.Execute("scratching", $true, $false, $false, $false,
                             $false, $true, 1, $false, "scraping", 2)
}
